$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '43.831.79'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '2.249.08'
$ws.Range('E3').Value = '  +0.24%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '322.51'
$ws.Range('E5').Value = '  +0.46%  '
$ws.Range('D6').Value = '101.78'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('E7').Value = '  -1.27%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -1.69%  '
$ws.Range('D10').Value = '37.11'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').Value = '7.72'
$ws.Range('E12').Value = '  -0.36%  '
$ws.Range('E13').Value = '  -2.57%  '
$ws.Range('D14').Value = '2.590.06'
$ws.Range('E14').Value = '  +0.53%  '
$ws.Range('D15').Value = '0.858'
$ws.Range('E15').Value = '  -1.50%  '
$ws.Range('B16').Value = 'Chainlink'
$ws.Range('C16').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D16').Value = '14.22'
$ws.Range('E16').Value = '  -1.82%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '2.248.68'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '43.740.81'
$ws.Range('E18').Value = '  +1.28%  '
$ws.Range('D19').Value = '13.63'
$ws.Range('E19').Value = '  -6.27%  '
$ws.Range('D20').Value = '0.0₃0989'
$ws.Range('E20').Value = '  +1.88%  '
$ws.Range('E21').Value = '  +0.03%  '
$ws.Range('D22').Value = '65.51'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('E23').Value = '  -1.10%  '
$ws.Range('D24').Value = '236.44'
$ws.Range('E24').Value = '  -0.91%  '
$ws.Range('E25').Value = '  -1.19%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').Value = '10.18'
$ws.Range('E27').Value = '  +0.84%  '
$ws.Range('E28').Value = '  -2.22%  '
$ws.Range('D29').Value = '37.17'
$ws.Range('E29').Value = '  +3.69%  '
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('D31').Value = '161.00'
$ws.Range('E31').Value = '  +4.54%  '
$ws.Range('D32').Value = '20.18'
$ws.Range('E32').Value = '  -1.73%  '
$ws.Range('D33').Value = '0.0852'
$ws.Range('E33').Value = '  -3.07%  '
$ws.Range('E34').Value = '  -1.12%  '
$ws.Range('D35').Value = '0.114'
$ws.Range('E35').Value = '  +9.39%  '
$ws.Range('D36').Value = '3.07'
$ws.Range('E36').Value = '  -3.44%  '
$ws.Range('D37').Value = '1.94'
$ws.Range('E37').Value = '  -2.05%  '
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('D39').Value = '3.77'
$ws.Range('E39').Value = '  +1.61%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '4.25'
$ws.Range('E40').Value = '  -4.99%  '
$ws.Range('B41').Value = 'Celestia'
$ws.Range('C41').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D41').Value = '16.01'
$ws.Range('E41').Value = '  +21.71%  '
$ws.Range('E42').Value = '  -2.41%  '
$ws.Range('D44').Value = '1.815.30'
$ws.Range('E44').Value = '  +1.82%  '
$ws.Range('D45').Value = '76.15'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = '0.200'
$ws.Range('E46').Value = '  -3.47%  '
$ws.Range('D47').Value = '82.53'
$ws.Range('E47').Value = '  -5.01%  '
$ws.Range('D48').Value = '5.22'
$ws.Range('E48').Value = '  -2.36%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '1.69'
$ws.Range('E49').Value = '  +4.86%  '
$ws.Range('B50').Value = 'MultiversX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D50').Value = '58.70'
$ws.Range('E50').Value = '  -1.32%  '
$ws.Range('D51').Value = '103.66'
$ws.Range('E51').Value = '  -0.32%  '
